$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Percent" column header
$ws.Range("F3").Value = "Percent"

# Percent-of-total formulas for each budget line (relative row ref, absolute total ref)
$ws.Range("F4").Formula = '=E4/$E$9'
$ws.Range("F5").Formula = '=E5/$E$9'
$ws.Range("F6").Formula = '=E6/$E$9'
$ws.Range("F7").Formula = '=E7/$E$9'
$ws.Range("F8").Formula = '=E8/$E$9'

# Grand total row: total of totals, and 100% check
$ws.Range("E9").Formula = '=B9+C9+D9'
$ws.Range("F9").Formula = '=E9/$E$9'

# Column width adjustments (A widened for labels, F sized for the new data)
$ws.Columns.Item(1).ColumnWidth = 11.436197916666666
$ws.Columns.Item(6).ColumnWidth = 9.983072916666666

# Move the active selection
$ws.Range("H8").Select()
